$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.2225601905
$ws.Range("C2").Value = 0.00000274522491327

$ws.Range("B3").Value = 11222.5601905
$ws.Range("C3").Value = 0.08235674739809999

$ws.Range("B4").Value = 41452.77057564985
$ws.Range("C4").Value = 0.3784580791558588

$ws.Range("B5").Value = 829.0554115129971
$ws.Range("C5").Value = 0.007569161583117176
